$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.051.28'
$ws.Range('E2').Value = '  -0.07%  '

$ws.Range('D3').Value = '2.402.83'
$ws.Range('E3').Value = '  -0.69%  '

$ws.Range('E4').Value = '  +0.54%  '

$ws.Range('D5').Value = '''567.37'
$ws.Range('E5').Value = '  -0.39%  '

$ws.Range('D6').Value = '''142.14'
$ws.Range('E6').Value = '  +1.97%  '

$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.34%  '

$ws.Range('D8').Value = '''0.539'
$ws.Range('E8').Value = '  +2.21%  '

$ws.Range('D9').Value = '2.410.89'
$ws.Range('E9').Value = '  +0.28%  '

$ws.Range('E10').Value = '  +1.81%  '

$ws.Range('E11').Value = '  -0.20%  '

$ws.Range('D12').Value = '''5.20'
$ws.Range('E12').Value = '  +2.54%  '

$ws.Range('D13').Value = '''0.346'
$ws.Range('E13').Value = '  +2.61%  '

$ws.Range('D14').Value = '''26.51'
$ws.Range('E14').Value = '  +1.61%  '

$ws.Range('E15').Value = '  -0.03%  '

$ws.Range('D16').Value = '2.841.00'
$ws.Range('E16').Value = '  -0.31%  '

$ws.Range('D17').Value = '60.854.79'
$ws.Range('E17').Value = '  -0.08%  '

$ws.Range('D18').Value = '2.415.27'
$ws.Range('E18').Value = '  +0.65%  '

$ws.Range('D19').Value = '''8.07'
$ws.Range('E19').Value = '  +2.84%  '

$ws.Range('D20').Value = '''10.70'
$ws.Range('E20').Value = '  +0.70%  '

$ws.Range('D21').Value = '''324.79'
$ws.Range('E21').Value = '  +0.50%  '

$ws.Range('E22').Value = '  +1.13%  '

$ws.Range('D23').Value = '''6.08'
$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '''1.00'
$ws.Range('E24').Value = '  -0.26%  '

$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').Value = '''1.92'
$ws.Range('E25').Value = '  +5.33%  '

$ws.Range('D26').Value = '''65.16'
$ws.Range('E26').Value = '  +0.73%  '

$ws.Range('D27').Value = '''589.65'
$ws.Range('E27').Value = '  +0.80%  '

$ws.Range('D28').Value = '''8.25'
$ws.Range('E28').Value = '  +0.14%  '

$ws.Range('D29').Value = '0.0₃0949'
$ws.Range('E29').Value = '  +2.12%  '

$ws.Range('E30').Value = '  -0.75%  '

$ws.Range('D31').Value = '''8.03'
$ws.Range('E31').Value = '  +2.33%  '

$ws.Range('E32').Value = '  +1.58%  '

$ws.Range('D33').Value = '''1.81'
$ws.Range('E33').Value = '  -0.34%  '

$ws.Range('D34').Value = '''0.134'
$ws.Range('E34').Value = '  +0.95%  '

$ws.Range('D35').Value = '''1.47'
$ws.Range('E35').Value = '  +4.16%  '

$ws.Range('D36').Value = '''0.998'
$ws.Range('E36').Value = '  -0.56%  '

$ws.Range('D37').Value = '''153.30'
$ws.Range('E37').Value = '  +1.11%  '

$ws.Range('E38').Value = '  +1.43%  '

$ws.Range('E39').Value = '  +0.40%  '

$ws.Range('E40').Value = '  +0.73%  '

$ws.Range('E41').Value = '  +1.88%  '

$ws.Range('E42').Value = '  -0.08%  '

$ws.Range('E43').Value = '  +1.09%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '''41.87'
$ws.Range('E44').Value = '  +1.82%  '

$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''2.49'
$ws.Range('E45').Value = '  +6.78%  '

$ws.Range('D46').Value = '0.0₆0284'
$ws.Range('E46').Value = '  +4.41%  '

$ws.Range('D47').Value = '''141.92'
$ws.Range('E47').Value = '  -0.39%  '

$ws.Range('D48').Value = '''3.53'
$ws.Range('E48').Value = '  +0.94%  '

$ws.Range('D49').Value = '''0.592'
$ws.Range('E49').Value = '  +0.90%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '''19.68'
$ws.Range('E50').Value = '  +1.28%  '

$ws.Range('B51').Value = 'Hedera'
$ws.Range('C51').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D51').Value = '''0.0510'
$ws.Range('E51').Value = '  +1.67%  '

